$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B: row 2 becomes a hard-coded value; rows 3:9 get a fresh formula fill ---
$ws.Range("B2").Value = 1.04
$ws.Range("B3").Formula = "=F3/M3"
$ws.Range("B4").Formula = "=F4/M4"
$ws.Range("B5").Formula = "=F5/M5"
$ws.Range("B6").Formula = "=F6/M6"
$ws.Range("B7").Formula = "=F7/M7"
$ws.Range("B8").Formula = "=F8/M8"
$ws.Range("B9").Formula = "=F9/M9"

# --- New Z-score header cells ---
$ws.Range("N1").Value = "Z1"
$ws.Range("K15").Value = " "
$ws.Range("O1").Value = "Z2"
$ws.Range("P1").Value = "Z3"
$ws.Range("Q1").Value = "Z4"

# Z1 header carries the "0.0000" numeric style already used by the N/O columns
$ws.Range("N1").NumberFormat = "0.0000"

# O1 inherited column O's old "0.0000" default style; restore it to the plain
# header look shared by the rest of row 1 (copy format from a plain header cell)
$ws.Range("M1").Copy()
$ws.Range("O1").PasteSpecial(-4122)

# --- New Z-score formulas: (value - AVE) / STDEV, one column per transducer ---
$ws.Range("N2").Formula = "=(F2-M2)/L2"
$ws.Range("O2").Formula = "=(G2-M2)/L2"
$ws.Range("P2").Formula = "=(H2-M2)/L2"
$ws.Range("Q2").Formula = "=(I2-M2)/L2"

$ws.Range("N3").Formula = "=(F3-M3)/L3"
$ws.Range("O3").Formula = "=(G3-M3)/L3"
$ws.Range("P3").Formula = "=(H3-M3)/L3"
$ws.Range("Q3").Formula = "=(I3-M3)/L3"

$ws.Range("N4").Formula = "=(F4-M4)/L4"
$ws.Range("O4").Formula = "=(G4-M4)/L4"
$ws.Range("P4").Formula = "=(H4-M4)/L4"
$ws.Range("Q4").Formula = "=(I4-M4)/L4"

$ws.Range("N5").Formula = "=(F5-M5)/L5"
$ws.Range("O5").Formula = "=(G5-M5)/L5"
$ws.Range("P5").Formula = "=(H5-M5)/L5"
$ws.Range("Q5").Formula = "=(I5-M5)/L5"

$ws.Range("N6").Formula = "=(F6-M6)/L6"
$ws.Range("O6").Formula = "=(G6-M6)/L6"
$ws.Range("P6").Formula = "=(H6-M6)/L6"
$ws.Range("Q6").Formula = "=(I6-M6)/L6"

$ws.Range("N7").Formula = "=(F7-M7)/L7"
$ws.Range("O7").Formula = "=(G7-M7)/L7"
$ws.Range("P7").Formula = "=(H7-M7)/L7"
$ws.Range("Q7").Formula = "=(I7-M7)/L7"

$ws.Range("N8").Formula = "=(F8-M8)/L8"
$ws.Range("O8").Formula = "=(G8-M8)/L8"
$ws.Range("P8").Formula = "=(H8-M8)/L8"
$ws.Range("Q8").Formula = "=(I8-M8)/L8"

$ws.Range("N9").Formula = "=(F9-M9)/L9"
$ws.Range("O9").Formula = "=(G9-M9)/L9"
$ws.Range("P9").Formula = "=(H9-M9)/L9"
$ws.Range("Q9").Formula = "=(I9-M9)/L9"

# Give the whole Z-score block (N2:Q9) the "0.0000" numeric format
$ws.Range("N2:Q9").NumberFormat = "0.0000"

# --- Final selection, matching where the user left the cursor ---
$ws.Range("K24").Select()
